$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New text for the "ProgramsTab" query cell (B2): the Website column is now
# derived via a CASE expression (program_acronym / program_link) instead of
# the plain prg.website column.
$newProgramQuery = @"
SELECT DISTINCT 
    prg.program_name AS "Program",
  CASE
    WHEN prg.program_link IS NOT NULL THEN prg.program_acronym
        ELSE prg.program_link
    END  AS "Website",
    prg.focus_area AS "Focus Area",
    prg.cancer_type AS "Cancer Type",
 CASE 
        WHEN prg.data_link IS NOT NULL THEN prg.website       
        ELSE prg.data_link
    END AS "Data Location Details"
FROM 
    df_program prg
WHERE 
     prg.cancer_type LIKE '%Lung Cancer%'
ORDER BY 
    lower(prg.program_name) ASC
LIMIT 100;
"@

# Touch the font (re-apply the theme color it already has) so that the cell
# picks up a distinct style record, matching the formatting refresh that was
# captured for this cell in the authored edit, then update its contents.
$ws.Range("B2").Font.ThemeColor = 1
$ws.Range("B2").Value = $newProgramQuery

# Update the saved selection/active cell to C3, scrolled back to the top of
# the sheet.
[void]$ws.Range("C3").Select()
